$p = $ppt.ActivePresentation

# --- Slide 3: Title placeholder "Couse Schedule" -> "Course Schedule" ---
# Original run has no dirty attr and is followed by a stray <a:endParaRPr>.
# Delete+InsertAfter rewrites the paragraph with a single run and drops the
# now-superfluous endParaRPr, matching the target markup shape.
$slide3 = $p.Slides.Item(3)
$title3 = $slide3.Shapes.Item(2)
$tr3 = $title3.TextFrame.TextRange
$tr3.Delete()
$tr3.InsertAfter("Course Schedule") | Out-Null

# --- Slide 4: Title placeholder "Couse Schedule" -> "Course Schedule" ---
$slide4 = $p.Slides.Item(4)
$title4 = $slide4.Shapes.Item(2)
$title4.TextFrame.TextRange.Text = "Course Schedule"

# --- Slide 5: Title placeholder "Couse Schedule" -> "Course Schedule" ---
$slide5 = $p.Slides.Item(5)
$title5 = $slide5.Shapes.Item(2)
$title5.TextFrame.TextRange.Text = "Course Schedule"

# --- Slide 5: Content placeholder "Weeks 8-10..." keeps its text but loses
# the stray trailing <a:endParaRPr>. Rewrite it the same way so the
# paragraph ends up with just the run.
$weeks5 = $slide5.Shapes.Item(3)
$trWeeks = $weeks5.TextFrame.TextRange
$weeksText = $trWeeks.Text
$trWeeks.Delete()
$trWeeks.InsertAfter($weeksText) | Out-Null
